# Add a new "2020" column (Q) to the tourism-GDP-share table, mirroring the
# formatting of the existing "2019" column (P), and update the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (year header): copy P4's format into Q4, then set its value ---
$yearSrc = $ws.Range("P4")
$yearDst = $ws.Range("Q4")
$yearSrc.Copy()
$yearDst.PasteSpecial(-4122)   # xlPasteFormats
$yearDst.Value = 2020

# --- Row 5 (percentage value): copy P5's format into Q5, then set its value ---
$valueSrc = $ws.Range("P5")
$valueDst = $ws.Range("Q5")
$valueSrc.Copy()
$valueDst.PasteSpecial(-4122)  # xlPasteFormats
$valueDst.Value = 3.3

# --- Update the active selection to reflect the new last-used cell ---
$ws.Range("R4").Select()
